$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.018.64'

$ws.Range("D3").Value = '2.268.00'
$ws.Range("E3").Value = '  +4.62%  '

$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.66%  '

$ws.Range("E6").Value = '  +1.72%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.31'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.87%  '

$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("E9").Value = '  +3.88%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '59.80'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.95%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0902'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.105'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.81%  '

$ws.Range("D13").Value = '2.608.05'
$ws.Range("E13").Value = '  +4.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.08'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.88'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.823'
$ws.Range("D16").Style = "Normal"

$ws.Range("E17").Value = '  +3.51%  '

$ws.Range("D18").Value = '2.270.31'
$ws.Range("E18").Value = '  +4.78%  '

$ws.Range("D19").Value = '41.924.28'
$ws.Range("E19").Value = '  +6.04%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '74.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.62%  '

$ws.Range("D21").Value = '0.0₃0932'
$ws.Range("E21").Value = '  +9.74%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.16'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.77%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.11'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +10.05%  '

$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.11%  '

$ws.Range("E27").Value = '  +7.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.86'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.11'
$ws.Range("D29").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.87%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.84'
$ws.Range("D32").Style = "Normal"

$ws.Range("E33").Value = '  +1.39%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.10'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.71%  '

$ws.Range("E35").Value = '  +4.26%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0640'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.09%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.37%  '

$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.81'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.62%  '

$ws.Range("E39").Value = '  +0.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.000255'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +44.60%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.13'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +17.04%  '

$ws.Range("E42").Value = '  +0.54%  '

$ws.Range("E43").Value = '  +6.01%  '

$ws.Range("E44").Value = '  +11.62%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.89'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.65%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.65%  '

$ws.Range("E47").Value = '  +3.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0983'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.24%  '

$ws.Range("D49").Value = '1.506.61'
$ws.Range("E49").Value = '  -1.38%  '

$ws.Range("E50").Value = '  +1.47%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.81'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.16%  '
